$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$installedApps = @'
[  'rest_framework',
    'simple_history',
    'django.contrib.admin',
    'django.contrib.auth',
    'django.contrib.contenttypes',
    'django.contrib.sessions',
    'django.contrib.messages',
    'django.contrib.staticfiles',
    'crispy_forms',
    'corsheaders',
    'debug_toolbar',
    'django_filters',
    'simple_history',
]
'@

$middleware = @'
[
    'django.middleware.security.SecurityMiddleware',
    'django.contrib.sessions.middleware.SessionMiddleware',
    'django.middleware.common.CommonMiddleware',
    'django.middleware.csrf.CsrfViewMiddleware',
    'django.contrib.auth.middleware.AuthenticationMiddleware',
    'django.contrib.messages.middleware.MessageMiddleware',
    'django.middleware.clickjacking.XFrameOptionsMiddleware',
    'simple_history.middleware.HistoryRequestMiddleware',]
'@

$ws.Range("B7").Value = $installedApps
$ws.Range("B8").Value = $middleware

$ws.Range("B10").Select()
